{"js": "// The document consists of a single title paragraph (a date, e.g.\n// \"2023-08-17 Thursday\") followed by a 20x5 table of cells, each\n// containing one paragraph with one arithmetic equation such as\n// \"45-8=37\". In document order, `context.document.body.paragraphs`\n// yields exactly: [title, cell(1,1), cell(1,2), ... cell(20,5)] \u2014 101\n// paragraphs total, matching the order the replacement text values\n// below were scraped from the canonical OOXML diff. We update each\n// paragraph's text positionally (not via global find/replace) because\n// a couple of the original equation strings repeat verbatim in\n// different cells but must map to different new values.\nconst newValues = [\n  \"2023-08-18 Friday\", \"79-9=70\", \"96-95=1\", \"86-67=19\", \"5+85=90\",\n  \"61-20=41\", \"7+65=72\", \"0+43=43\", \"16+39=55\", \"56+19=75\",\n  \"39-23=16\", \"53-2=51\", \"46-24=22\", \"73+4=77\", \"91-12=79\",\n  \"20+60=80\", \"39+44=83\", \"6+82=88\", \"32+20=52\", \"48+46=94\",\n  \"69-52=17\", \"41+31=72\", \"73-72=1\", \"34+52=86\", \"48+8=56\",\n  \"61+27=88\", \"51-45=6\", \"3+29=32\", \"84+8=92\", \"84+10=94\",\n  \"77+4=81\", \"25+3=28\", \"91-49=42\", \"86+8=94\", \"9+54=63\",\n  \"43+10=53\", \"95-28=67\", \"33-15=18\", \"36+55=91\", \"80-71=9\",\n  \"75-68=7\", \"9+1=10\", \"89-61=28\", \"42-34=8\", \"71+24=95\",\n  \"67-13=54\", \"51+34=85\", \"28+25=53\", \"88-29=59\", \"53-53=0\",\n  \"51-33=18\", \"1+50=51\", \"32+21=53\", \"60-43=17\", \"90-16=74\",\n  \"91-42=49\", \"30-7=23\", \"98-69=29\", \"2+57=59\", \"90-52=38\",\n  \"28+55=83\", \"55-7=48\", \"54-33=21\", \"97-27=70\", \"67-6=61\",\n  \"35+6=41\", \"56-45=11\", \"2+27=29\", \"75-68=7\", \"5+55=60\",\n  \"61+0=61\", \"43+14=57\", \"50-12=38\", \"24+21=45\", \"17+39=56\",\n  \"86-39=47\", \"30+22=52\", \"14+74=88\", \"12-4=8\", \"78-17=61\",\n  \"2+73=75\", \"74+13=87\", \"82-54=28\", \"24-3=21\", \"68-1=67\",\n  \"96-51=45\", \"56+15=71\", \"16-0=16\", \"73-38=35\", \"84-71=13\",\n  \"18-4=14\", \"7+82=89\", \"97+2=99\", \"4-0=4\", \"63-25=38\",\n  \"69+23=92\", \"0+31=31\", \"1+71=72\", \"89-18=71\", \"6+53=59\",\n  \"63+34=97\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newValues.length +\n    \" but found \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < newValues.length; i++) {\n  paragraphs.items[i].insertText(newValues[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document is a single title paragraph (a date, e.g.\n# \"2023-08-17 Thursday\") followed by one 20-row x 5-column table whose\n# cells each hold one arithmetic equation such as \"45-8=37\". We update\n# the title text and then walk the table cell-by-cell (row-major,\n# matching the canonical OOXML order) so that each cell gets its own\n# new value positionally -- a plain global find/replace would be\n# unsafe here because a couple of the original equations (e.g.\n# \"3+77=80\") occur more than once but must map to different results.\n\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Item(1).Range.Text = \"2023-08-18 Friday\"\n\n$newValues = @(\n    @(\"79-9=70\", \"96-95=1\", \"86-67=19\", \"5+85=90\", \"61-20=41\"),\n    @(\"7+65=72\", \"0+43=43\", \"16+39=55\", \"56+19=75\", \"39-23=16\"),\n    @(\"53-2=51\", \"46-24=22\", \"73+4=77\", \"91-12=79\", \"20+60=80\"),\n    @(\"39+44=83\", \"6+82=88\", \"32+20=52\", \"48+46=94\", \"69-52=17\"),\n    @(\"41+31=72\", \"73-72=1\", \"34+52=86\", \"48+8=56\", \"61+27=88\"),\n    @(\"51-45=6\", \"3+29=32\", \"84+8=92\", \"84+10=94\", \"77+4=81\"),\n    @(\"25+3=28\", \"91-49=42\", \"86+8=94\", \"9+54=63\", \"43+10=53\"),\n    @(\"95-28=67\", \"33-15=18\", \"36+55=91\", \"80-71=9\", \"75-68=7\"),\n    @(\"9+1=10\", \"89-61=28\", \"42-34=8\", \"71+24=95\", \"67-13=54\"),\n    @(\"51+34=85\", \"28+25=53\", \"88-29=59\", \"53-53=0\", \"51-33=18\"),\n    @(\"1+50=51\", \"32+21=53\", \"60-43=17\", \"90-16=74\", \"91-42=49\"),\n    @(\"30-7=23\", \"98-69=29\", \"2+57=59\", \"90-52=38\", \"28+55=83\"),\n    @(\"55-7=48\", \"54-33=21\", \"97-27=70\", \"67-6=61\", \"35+6=41\"),\n    @(\"56-45=11\", \"2+27=29\", \"75-68=7\", \"5+55=60\", \"61+0=61\"),\n    @(\"43+14=57\", \"50-12=38\", \"24+21=45\", \"17+39=56\", \"86-39=47\"),\n    @(\"30+22=52\", \"14+74=88\", \"12-4=8\", \"78-17=61\", \"2+73=75\"),\n    @(\"74+13=87\", \"82-54=28\", \"24-3=21\", \"68-1=67\", \"96-51=45\"),\n    @(\"56+15=71\", \"16-0=16\", \"73-38=35\", \"84-71=13\", \"18-4=14\"),\n    @(\"7+82=89\", \"97+2=99\", \"4-0=4\", \"63-25=38\", \"69+23=92\"),\n    @(\"0+31=31\", \"1+71=72\", \"89-18=71\", \"6+53=59\", \"63+34=97\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
